$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 24: fill in the Finance and Category Screen / Yoddle question
$ws.Range("B24").Value = "Finance and Category Screen"
$ws.Range("C24").Value = "Yoddle"
$ws.Range("F24").Value = "What all endpoints we will be using for Income and expense for a category?"

# Row 25: fill in the Category Screen / Yoddle question
$ws.Range("B25").Value = "Category Secreen"
$ws.Range("C25").Value = "Yoddle"
$ws.Range("F25").Value = "What all endpoints we will be using for first level category on category Screen"

# Row 26 (new row): Transaction / Yoddle question
$ws.Range("A26").Value = "PFM"
$ws.Range("B26").Value = "Transaction"
$ws.Range("C26").Value = "Yoddle"
$ws.Range("F26").Value = "What all endpoints we will be using for second level category on category Screen"

# Widen column F to fit the new, longer text
$ws.Columns("F").ColumnWidth = 76.8

# Update the view: scroll so row 19 is at the top, and select B26 (last edited cell)
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
$ws.Range("B26").Select()
